# Update BatchInput worksheet to showcase all variations currently possible.
# Rows 2-6 get new parameter values and rows 7-11 are newly added, sweeping
# the xRep parameter (column B) from 5 through 10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.03996003996003996
$ws.Range("G2").Value = 0.03996003996003996
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 0.03996003996003996
$ws.Range("G3").Value = 0.03996003996003996
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.03996003996003996
$ws.Range("G4").Value = 0.03996003996003996
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.5
$ws.Range("J4").Value = 0

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.03996003996003996
$ws.Range("G5").Value = 0.03996003996003996
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 0

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.03996003996003996
$ws.Range("G6").Value = 0.03996003996003996
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0.5

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.03996003996003996
$ws.Range("G7").Value = 0.03996003996003996
$ws.Range("H7").Value = 0.03996003996003996
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 0

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.04595404595404595
$ws.Range("G8").Value = 0.04595404595404595
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 0

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.04595404595404595
$ws.Range("G9").Value = 0.04595404595404595
$ws.Range("H9").Value = 0.04595404595404595
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.04595404595404595
$ws.Range("G10").Value = 0.04595404595404595
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.07992007992007992
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0

# Apply the "0.000" number format to the F:H columns across all data rows
# (reuses the existing style already present on F2:H6 in the template).
$ws.Range("F2:H11").NumberFormat = "0.000"

# I2:J3 and the new empty K6 cell also pick up the "0.000" style in the target.
$ws.Range("I2:J3").NumberFormat = "0.000"
$ws.Range("K6").NumberFormat = "0.000"

# Restore the selection to match the edited region.
$ws.Range("F10:H10").Select()

